# Repull data, push all data, mean calculation
# Update the dSF column (F) values to reflect recalculated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 4
    5  = 4
    6  = 3
    7  = 2
    8  = -3
    9  = 1
    10 = 5
    11 = -2
    12 = 2
    13 = 1
    14 = -1
    15 = 3
    16 = 3
    17 = 4
    18 = -1
    19 = 10
    20 = 2
    21 = 2
    22 = 4
    23 = 1
    24 = 4
    25 = -5
    26 = -2
    27 = -2
    28 = 3
    29 = -2
    30 = -2
    31 = -1
    33 = 1
    34 = 3
    35 = -8
    37 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
